# Edit script: remove the (now redundant) raw-numeric "Cout optimal" column,
# shift the remaining data left, fix up column order, drop the trailing
# empty rows, and re-point the selection - matching the authoring app's
# regenerated results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the old column B ("Cout optimal" raw numbers) entirely; this
#    shifts old C -> B, old D -> C, old E -> D.
$ws.Range("B1").EntireColumn.Delete()

# 2) The regenerated sheet actually orders the remaining two data columns
#    the other way round (Cout algo final before Cout optimal), so swap
#    the content of columns C and D for every data row (1-10).
for ($r = 1; $r -le 10; $r++) {
    $c = $ws.Cells.Item($r, 3)
    $d = $ws.Cells.Item($r, 4)
    $cVal = $c.Value2
    $dVal = $d.Value2
    $c.Value2 = $dVal
    $d.Value2 = $cVal
}

# 3) Row 8's middle column picks up a freshly computed text value rather
#    than reusing the "maison manquante..." label.
$ws.Cells.Item(8, 3).Value2 = "5.094179894179893"

# 4) Drop the long tail of empty rows (11-33): only 10 data rows remain.
$ws.Range("A11:E33").EntireRow.Delete()

# 5) Restore the default top-left cell / selection as in the saved file.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
